$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 2.15
$ws.Range("L4").Value = 2.7
$ws.Range("Y4").Value = 10.75
$ws.Range("Z4").Value = 40
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 5.1
$ws.Range("AO4").Value = 16
$ws.Range("AR4").Value = 100
$ws.Range("AX4").Value = 10.75

$wb.Save()
